$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 58, pushing the existing rows 58..133 down to 59..134.
# Excel's Insert() copies the formatting (e.g. the date number-format on column D)
# from the row above, matching the style already used throughout the table.
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new weekly price observation.
$ws.Cells.Item(58, 1).Value2  = 1
$ws.Cells.Item(58, 2).Value2  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(58, 3).Value2  = 'Arica y Parinacota'
$ws.Cells.Item(58, 4).Value2  = 44757
$ws.Cells.Item(58, 5).Value2  = 15
$ws.Cells.Item(58, 6).Value2  = 'Fruta'
$ws.Cells.Item(58, 7).Value2  = 100108
$ws.Cells.Item(58, 8).Value2  = 'Tropicales y subtropicales'
$ws.Cells.Item(58, 9).Value2  = 100108003
$ws.Cells.Item(58, 10).Value2 = 'Maracuyá'
$ws.Cells.Item(58, 11).Value2 = 'Sin especificar'
$ws.Cells.Item(58, 12).Value2 = 'Primera'
$ws.Cells.Item(58, 13).Value2 = 120
$ws.Cells.Item(58, 14).Value2 = 22000
$ws.Cells.Item(58, 15).Value2 = 23000
$ws.Cells.Item(58, 16).Value2 = 22500
$ws.Cells.Item(58, 17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(58, 18).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(58, 19).Value2 = 1125
$ws.Cells.Item(58, 20).Value2 = 20
